$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bold the header row (A1:D1, L1) ---
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("L1").Font.Bold = $true

# --- Shift the "Features" list in column L up (rows 5-8), clearing the old tail ---
$ws.Range("L5").Value = "Powerups"
$ws.Range("L6").Value = "Obstacles"
$ws.Range("L7").Value = "Main Menu"
$ws.Range("L8").Value = "Score"
$ws.Range("L12:L14").ClearContents()

# --- Rework the "Powerups" block (was A13:B14, now A12:C14 with extra data) ---
$ws.Range("A13").ClearContents()
$ws.Range("A12").Value = "Powerups"
$ws.Range("B12").Value = "Life Pickup"
$ws.Range("C12").Value = 1
$ws.Range("B13").Value = "Chase Pickup"
$ws.Range("C13").Value = "isChase=true"
$ws.Range("B14").Value = "Point Pickup"
$ws.Range("C14").Value = 1

# --- New "Scoring" block ---
$ws.Range("A15").Value = "Scoring"
$ws.Range("B15").Value = "Distance"
$ws.Range("B16").Value = "Points"

# --- New "State Machine" block ---
$ws.Range("A17").Value = "State Machine"
$ws.Range("B17").Value = "Main Menu state"
$ws.Range("B18").Value = "Gameplay state"
$ws.Range("B19").Value = "Pause state"
$ws.Range("B20").Value = "Defeated state"

# --- New "Chaser" block ---
$ws.Range("A21").Value = "Chaser"
$ws.Range("B21").Value = "Number of spawns"
$ws.Range("B22").Value = "Spawn Rate"
$ws.Range("B23").Value = "Damage"
$ws.Range("B24").Value = "Speed"

# --- New "Obstacles" block ---
$ws.Range("A25").Value = "Obstacles"
$ws.Range("B25").Value = "Number of spawns"
$ws.Range("B26").Value = "Spawn Rate"
$ws.Range("B27").Value = "Damage"
$ws.Range("B28").Value = "Speed"

# --- Trailing notes ---
$ws.Range("A29").Value = "Spawn Lanes"
$ws.Range("A30").Value = "Chasing mode"

# --- Restore the active cell selection ---
[void]$ws.Range("E16").Select()
